$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 10 corresponds to "Mostrar historial de asistencias de un grupo." task
# Update status from "Por iniciar" to "Hecho"
$ws.Range("F10").Value = "Hecho"

# Register 4 consumed hours on "Dia 12" column (AC) for row 10
$ws.Range("AC10").Value = 4

# Move active selection to AC11 as in the final workbook
$ws.Range("AC11").Select()
